$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ExploreArea": add an "Unlock" column and a new AreaID 101 row.
# ---------------------------------------------------------------------------
$wsArea = $wb.Worksheets.Item("ExploreArea")

# Insert a new column D (shifts the old "ExploreList" column from D to E)
$wsArea.Columns.Item(4).Insert() | Out-Null
$wsArea.Columns.Item(4).ColumnWidth = 16.8

# Header row
$wsArea.Cells.Item(1,4).Value = "Unlock"

# Row 2 (AreaID 100) now unlocked by default
$wsArea.Cells.Item(2,4).Value = $true

# Row 3: brand new explore area
$wsArea.Cells.Item(3,1).Value = 101
$wsArea.Cells.Item(3,2).Value = "ExploreArea_Name_101"
$wsArea.Cells.Item(3,3).Value = "ExploreArea_Desc_101"
$wsArea.Cells.Item(3,4).Value = $false
$wsArea.Cells.Item(3,5).Value = "1,2,3,4,5"

# ---------------------------------------------------------------------------
# Restore the selections / active sheet recorded by the author while editing.
# ---------------------------------------------------------------------------
$wsArea.Activate() | Out-Null
$wsArea.Range("C14").Select() | Out-Null

$wsData = $wb.Worksheets.Item("ExploreData")
$wsData.Activate() | Out-Null
$wsData.Range("B29").Select() | Out-Null

$wsPoint = $wb.Worksheets.Item("ExplorePoint")
$wsPoint.Activate() | Out-Null
$wsPoint.Range("D2:D15").Select() | Out-Null

$wsEvent = $wb.Worksheets.Item("ExploreEvent")
$wsEvent.Activate() | Out-Null
$wsEvent.Range("B6").Select() | Out-Null

$wsChoose = $wb.Worksheets.Item("ExploreChoose")
$wsChoose.Activate() | Out-Null
$wsChoose.Range("B7").Select() | Out-Null

# Leave the ExploreArea sheet as the active tab, matching the target state.
$wsArea.Activate() | Out-Null
$wsArea.Range("C14").Select() | Out-Null
